$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Companies")

$ws.Range("A25").Value = "CAREFIRST"
$ws.Range("A28").Value = "EXCELLUS"
$ws.Range("A29").Value = "HARMARK"
$ws.Range("A31").Value = "HIGHMARK"
$ws.Range("A32").Value = "HORIZON"
$ws.Range("A34").Value = "PREMERA"
$ws.Range("A36").Value = "WELLMARK"
$ws.Range("A2").Value = "ANTHEM"
